# Update the "datetimeFigureOut" date field placeholders (slide master, all
# slide layouts, and the notes master) from the stale "7/20/17" to "4/16/2018".
$p = $ppt.ActivePresentation
$OLD_DATE = "7/20/17"
$NEW_DATE = "4/16/2018"

function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $OLD_DATE) {
                $sh.TextFrame.TextRange.Text = $NEW_DATE
            }
            # Some containers (e.g. the notes master) only accept the edit
            # through the linked header/footer date field rather than by
            # writing the placeholder's TextRange directly, so fall back to
            # that API whenever the direct edit above did not stick.
            if ($sh.TextFrame.TextRange.Text -eq $OLD_DATE) {
                $container.HeadersFooters.DateAndTime.Text = $NEW_DATE
            }
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

# Slide master date placeholder.
Update-DatePlaceholder $master

# Every slide layout's date placeholder.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout
}

# Notes master date placeholder.
Update-DatePlaceholder $p.NotesMaster

# Remove the obsolete "UndoRedo Stack" diagram (rectangle + connector +
# annotation textbox) from the single content slide; this portion of the
# Logic component diagram documented an UndoRedoStack design that no longer
# exists now that undo/redo is handled by VersionedAddressBook.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "TextBox 62" -or $sh.Name -eq "Straight Arrow Connector 57" -or $sh.Name -eq "Rectangle 62") {
        if ($sh.Id -eq 63 -or $sh.Id -eq 61 -or $sh.Id -eq 59) {
            $sh.Delete()
        }
    }
}
